$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the same formatting as the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Threshold observed for a save: the recurring "sum" value in column G
$saveThreshold = 8.418600821238126

# Populate H2:H61 with 1 if the "sum" (column G) meets/exceeds the save
# threshold, else 0.
for ($r = 2; $r -le 61; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -ge $saveThreshold) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
